$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting N:P -> O:Q.
$ws.Columns("N").Insert()

# The newly inserted column should take on a width close to column M's
# (10.7109375 chars), matching Excel's default "insert column" behaviour
# of carrying over the format/width of the column to the left.
$ws.Columns("N").ColumnWidth = 9.83

# Make "Repayment schedule" the active/selected sheet and update its
# selected cell.
$ws.Activate()
$ws.Range("J18").Select() | Out-Null
